$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = 43809.22349537037
$ws.Range("B42").Value = "jetnew"
$ws.Range("C42").Value = "test"
$ws.Range("D42").Value = "No content found."

$ws.Range("A43").Value = 43809.22594907408
$ws.Range("B43").Value = "jetnew"
$ws.Range("C43").Value = "test"
$ws.Range("D43").Value = "No content found."

$ws.Range("A44").Value = 43809.22703703704
$ws.Range("B44").Value = "jetnew"
$ws.Range("C44").Value = "how long do students live in cinnamon college?"
$ws.Range("D44").Value = "two years"

$ws.Range("A42:A44").NumberFormat = "yyyy-mm-dd h:mm:ss"
